$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Valor Mora" values between the row for period 2402 (row 17)
# and the row for period 2308 (row 23) as part of updating the EC database.
$ws.Range("F17").Value = 31280
$ws.Range("F23").Value = 18768
